# Apply "Share"/"Total" columns to each GRA-* sheet, and update the
# GRA-carbontax / GRA-fueltax weight values per the authored commit.

$wb = $excel.ActiveWorkbook

# Sheet name -> fixed "Total" value used for every row of that sheet.
$totals = @{
    "GRA-carbontax"           = 15
    "GRA-fueltax"              = 7.5
    "GRA-evsubsidy"            = 5
    "GRA-elecgensubsidy"       = 5
    "GRA-eleccapconstsubsidy"  = 5
    "GRA-distsolarsubsidy"     = 5
    "GRA-fuelsubsidy"          = 5
    "GRA-ntnldebtinterest"     = 10
    "GRA-remainder"            = 15
}

# Sheets whose new "Total"/"Share" header cells (C1/D1) pick up the
# darker/black header font used elsewhere in the workbook.
$darkHeaderSheets = @(
    "GRA-evsubsidy",
    "GRA-elecgensubsidy",
    "GRA-eleccapconstsubsidy",
    "GRA-distsolarsubsidy",
    "GRA-fuelsubsidy",
    "GRA-ntnldebtinterest",
    "GRA-remainder"
)

$sheetOrder = @(
    "GRA-carbontax",
    "GRA-fueltax",
    "GRA-evsubsidy",
    "GRA-elecgensubsidy",
    "GRA-eleccapconstsubsidy",
    "GRA-distsolarsubsidy",
    "GRA-fuelsubsidy",
    "GRA-ntnldebtinterest",
    "GRA-remainder"
)

foreach ($name in $sheetOrder) {
    $ws = $wb.Worksheets.Item($name)

    # Headers (order matters for shared-string table ordering: Share then Total)
    $ws.Cells.Item(1,4).Value = "Share"
    $ws.Cells.Item(1,3).Value = "Total"

    if ($darkHeaderSheets -contains $name) {
        $ws.Range("C1:D1").Font.Color = 0
    }

    $total = $totals[$name]

    for ($r = 2; $r -le 6; $r++) {
        $ws.Cells.Item($r,3).Value = $total
        $ws.Cells.Item($r,4).Formula = "=B$r/C$r"
        $ws.Cells.Item($r,4).NumberFormat = "0.00%"
    }
}

# GRA-carbontax: replace the TRANSPOSE array formula with literal weights
$wsCarbon = $wb.Worksheets.Item("GRA-carbontax")
$carbonVals = @(5, 2.5, 10, 2.5, 0)
for ($i = 0; $i -lt 5; $i++) {
    $r = 2 + $i
    $wsCarbon.Cells.Item($r,2).Value = $carbonVals[$i]
    $wsCarbon.Cells.Item($r,2).Font.Name = "Helvetica Neue"
    $wsCarbon.Cells.Item($r,2).Font.Size = 10
    $wsCarbon.Cells.Item($r,2).Font.Color = 0
}

# GRA-fueltax: replace the TRANSPOSE array formula with literal weights
$wsFuel = $wb.Worksheets.Item("GRA-fueltax")
$fuelVals = @(0, 2.5, 0, 0, 5)
for ($i = 0; $i -lt 5; $i++) {
    $r = 2 + $i
    $wsFuel.Cells.Item($r,2).Value = $fuelVals[$i]
}
